$wb = $excel.ActiveWorkbook

# Work on the "RangeTests" sheet (second sheet / tab) which holds the
# ChanSlope data.
$ws = $wb.Worksheets.Item("RangeTests")
$ws.Activate()

# Update the slope exponent in A1 - this drives every formula in the grid.
$ws.Range("A1").Value = -0.05

# Move the selection to A2 to match the new cursor position.
$ws.Range("A2").Select()

$wb.RecalculateFull()
